$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row to append (row 22), mirroring the existing table's text-formatted columns.
# Columns A (date-looking) and C (numeric-looking) must be forced to Text format
# first so Excel doesn't silently coerce them into a date serial / number.
$ws.Cells.Item(22, 1).NumberFormat = "@"
$ws.Cells.Item(22, 3).NumberFormat = "@"

$ws.Cells.Item(22, 1).Value = "2025-10-08"
$ws.Cells.Item(22, 2).Value = "Pick 4"
$ws.Cells.Item(22, 3).Value = "251008"
$ws.Cells.Item(22, 4).Value = "2-0-8-2"
$ws.Cells.Item(22, 5).Value = "2025-10-08T21:38:15.767+04:00"
